# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2026-02-07 (serial 46060) to 2026-02-08 (serial 46061).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on column A (Beteckning), data starts at row 2.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
